$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.TopLeftCell = "A53"
Write-Output ("Readback TopLeftCell: " + $excel.ActiveWindow.TopLeftCell)
$ws.Range("Q73").Select()
